$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a single-space value into C3 (new shared string), which extends the used range to A1:C3
$ws.Range("C3").Value = " "

# Move the active selection to D3 (mirrors the author's next click after filling C3)
$ws.Range("D3").Select()
